# Updates TPM-derived metrics (columns G-J, M-T) for rows 2-10
# to reflect the new TPM normalization values referenced in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.120168333333333
$ws.Range("H2").Value = 3.360505
$ws.Range("I2").Value = 0.001768092629909379
$ws.Range("J2").Value = 0.001768092629909379
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 86.3532428697611
$ws.Range("R2").Value = 777.17918582785
$ws.Range("S2").Value = 0.0004250170544469675
$ws.Range("T2").Value = 0.0004250170544469675
# Row 3
$ws.Range("G3").Value = 1.120168333333333
$ws.Range("H3").Value = 3.360505
$ws.Range("I3").Value = 0.001768092629909379
$ws.Range("J3").Value = 0.001768092629909379
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 113.7867411196178
$ws.Range("R3").Value = 1024.08067007656
$ws.Range("S3").Value = 0.0005600404100482784
$ws.Range("T3").Value = 0.0005600404100482784
# Row 4
$ws.Range("G4").Value = 1.120168333333333
$ws.Range("H4").Value = 3.360505
$ws.Range("I4").Value = 0.001768092629909379
$ws.Range("J4").Value = 0.001768092629909379
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 159.0939118962044
$ws.Range("R4").Value = 1431.84520706584
$ws.Range("S4").Value = 0.0007830351654141328
$ws.Range("T4").Value = 0.0007830351654141328
# Row 5
$ws.Range("I5").Value = 0.9534130698726969
$ws.Range("J5").Value = 0.9534130698726969
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 46564.47800596355
$ws.Range("R5").Value = 419080.302053672
$ws.Range("S5").Value = 0.2291830234309066
$ws.Range("T5").Value = 0.2291830234309066
# Row 6
$ws.Range("I6").Value = 0.9534130698726969
$ws.Range("J6").Value = 0.9534130698726969
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.3019920096744366
$ws.Range("T6").Value = 0.3019920096744366
# Row 7
$ws.Range("I7").Value = 0.9534130698726969
$ws.Range("J7").Value = 0.9534130698726969
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.4222380367673537
$ws.Range("T7").Value = 0.4222380367673537
# Row 8
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04481883749739363
$ws.Range("J8").Value = 0.04481883749739363
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 2188.941854110411
$ws.Range("R8").Value = 19700.4766869937
$ws.Range("S8").Value = 0.01077362688732878
$ws.Range("T8").Value = 0.01077362688732878
# Row 9
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04481883749739363
$ws.Range("J9").Value = 0.04481883749739363
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("Q9").Value = 2884.345182672658
$ws.Range("S9").Value = 0.01419629249357483
$ws.Range("T9").Value = 0.01419629249357482
# Row 10
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04481883749739363
$ws.Range("J10").Value = 0.04481883749739363
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("R10").Value = 36295.40476066288
$ws.Range("S10").Value = 0.01984891811649002
$ws.Range("T10").Value = 0.01984891811649002
